$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current ("before") values for rows 6-10, columns A,B,D,E,F,G,H,Q,R
$rows = 6,7,8,9,10
$cols = 'A','B','D','E','F','G','H','Q','R'

$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{}
    foreach ($c in $cols) {
        $before[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# The data for rows 6-10 cyclically shifts up by one row:
# new row 6  <- old row 7
# new row 7  <- old row 8
# new row 8  <- old row 9
# new row 9  <- old row 10
# new row 10 <- old row 6
$mapping = @{ 6 = 7; 7 = 8; 8 = 9; 9 = 10; 10 = 6 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $before[$src][$c]
    }
}
